$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp on the last existing row (row 44, column A)
$ws.Cells.Item(44, 1).Value = 44357.77969643634

# Append the new data row (row 45)
$ws.Cells.Item(45, 1).Value = 44358.76795423555
$ws.Cells.Item(45, 2).Value = 77078
$ws.Cells.Item(45, 3).Value = 64817
$ws.Cells.Item(45, 4).Value = 3344
$ws.Cells.Item(45, 5).Value = 2088
$ws.Cells.Item(45, 6).Value = 1462
$ws.Cells.Item(45, 7).Value = 20335
$ws.Cells.Item(45, 8).Value = 1406
$ws.Cells.Item(45, 9).Value = 866
$ws.Cells.Item(45, 10).Value = 177

# Match the date formatting style used in column A (style index 2 -> numFmt 164)
$ws.Cells.Item(45, 1).NumberFormat = $ws.Cells.Item(44, 1).NumberFormat
